$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Ligand symbol = Hspg2) and C (Receptor symbol = Col13a1) are constant
# across all data rows; fill them for the full new range first.
$ws.Range("B2:B11").Value = "Hspg2"
$ws.Range("C2:C11").Value = "Col13a1"

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 230.4028776666667
$ws.Range("H2").Value = 691.208633
$ws.Range("I2").Value = 0.5806109522726741
$ws.Range("J2").Value = 0.580610952272674
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.335983
$ws.Range("N2").Value = 1.007949
$ws.Range("O2").Value = 0.4945457382278176
$ws.Range("P2").Value = 0.4945457382278176
$ws.Range("Q2").Value = 77.41145004707965
$ws.Range("R2").Value = 696.703050423717
$ws.Range("S2").Value = 0.2871386720148458
$ws.Range("T2").Value = 0.2871386720148457

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 230.4028776666667
$ws.Range("H3").Value = 691.208633
$ws.Range("I3").Value = 0.5806109522726741
$ws.Range("J3").Value = 0.580610952272674
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.343394
$ws.Range("N3").Value = 1.030182
$ws.Range("O3").Value = 0.5054542617721824
$ws.Range("P3").Value = 0.5054542617721824
$ws.Range("Q3").Value = 79.11896577346732
$ws.Range("R3").Value = 712.070691961206
$ws.Range("S3").Value = 0.2934722802578283
$ws.Range("T3").Value = 0.2934722802578283

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 134.7127306666667
$ws.Range("H4").Value = 404.138192
$ws.Range("I4").Value = 0.339473567464654
$ws.Range("J4").Value = 0.339473567464654
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.335983
$ws.Range("N4").Value = 1.007949
$ws.Range("O4").Value = 0.4945457382278176
$ws.Range("P4").Value = 0.4945457382278176
$ws.Range("Q4").Value = 45.26118738757866
$ws.Range("R4").Value = 407.350686488208
$ws.Range("S4").Value = 0.1678852060306382
$ws.Range("T4").Value = 0.1678852060306381

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 134.7127306666667
$ws.Range("H5").Value = 404.138192
$ws.Range("I5").Value = 0.339473567464654
$ws.Range("J5").Value = 0.339473567464654
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.343394
$ws.Range("N5").Value = 1.030182
$ws.Range("O5").Value = 0.5054542617721824
$ws.Range("P5").Value = 0.5054542617721824
$ws.Range("Q5").Value = 46.25954343454933
$ws.Range("R5").Value = 416.335890910944
$ws.Range("S5").Value = 0.1715883614340159
$ws.Range("T5").Value = 0.1715883614340159

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1219473333333333
$ws.Range("H6").Value = 0.365842
$ws.Range("I6").Value = 0.0003073050043941503
$ws.Range("J6").Value = 0.0003073050043941503
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.335983
$ws.Range("N6").Value = 1.007949
$ws.Range("O6").Value = 0.4945457382278176
$ws.Range("P6").Value = 0.4945457382278176
$ws.Range("Q6").Value = 0.04097223089533333
$ws.Range("R6").Value = 0.368750078058
$ws.Range("S6").Value = 0.0001519763802592078
$ws.Range("T6").Value = 0.0001519763802592078

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1219473333333333
$ws.Range("H7").Value = 0.365842
$ws.Range("I7").Value = 0.0003073050043941503
$ws.Range("J7").Value = 0.0003073050043941503
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.343394
$ws.Range("N7").Value = 1.030182
$ws.Range("O7").Value = 0.5054542617721824
$ws.Range("P7").Value = 0.5054542617721824
$ws.Range("Q7").Value = 0.04187598258266666
$ws.Range("R7").Value = 0.376883843244
$ws.Range("S7").Value = 0.0001553286241349425
$ws.Range("T7").Value = 0.0001553286241349425

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1193106666666667
$ws.Range("H8").Value = 0.357932
$ws.Range("I8").Value = 0.0003006606535958338
$ws.Range("J8").Value = 0.0003006606535958338
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.335983
$ws.Range("N8").Value = 1.007949
$ws.Range("O8").Value = 0.4945457382278176
$ws.Range("P8").Value = 0.4945457382278176
$ws.Range("Q8").Value = 0.04008635571866666
$ws.Range("R8").Value = 0.360777201468
$ws.Range("S8").Value = 0.0001486904448886097
$ws.Range("T8").Value = 0.0001486904448886097

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1193106666666667
$ws.Range("H9").Value = 0.357932
$ws.Range("I9").Value = 0.0003006606535958338
$ws.Range("J9").Value = 0.0003006606535958338
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.343394
$ws.Range("N9").Value = 1.030182
$ws.Range("O9").Value = 0.5054542617721824
$ws.Range("P9").Value = 0.5054542617721824
$ws.Range("Q9").Value = 0.04097056706933334
$ws.Range("R9").Value = 0.368735103624
$ws.Range("S9").Value = 0.000151970208707224
$ws.Range("T9").Value = 0.000151970208707224

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 31.471469
$ws.Range("H10").Value = 94.414407
$ws.Range("I10").Value = 0.07930751460468206
$ws.Range("J10").Value = 0.07930751460468205
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.335983
$ws.Range("N10").Value = 1.007949
$ws.Range("O10").Value = 0.4945457382278176
$ws.Range("P10").Value = 0.4945457382278176
$ws.Range("Q10").Value = 10.573878569027
$ws.Range("R10").Value = 95.16490712124299
$ws.Range("S10").Value = 0.03922119335718592
$ws.Range("T10").Value = 0.03922119335718591

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 31.471469
$ws.Range("H11").Value = 94.414407
$ws.Range("I11").Value = 0.07930751460468206
$ws.Range("J11").Value = 0.07930751460468205
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.343394
$ws.Range("N11").Value = 1.030182
$ws.Range("O11").Value = 0.5054542617721824
$ws.Range("P11").Value = 0.5054542617721824
$ws.Range("Q11").Value = 10.807113625786
$ws.Range("R11").Value = 97.26402263207399
$ws.Range("S11").Value = 0.04008632124749615
$ws.Range("T11").Value = 0.04008632124749614
